$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

# The sheet currently has 9 rows: row 2 holds 'num_credito' / '805117974',
# and rows 3-9 hold the remaining stat rows (entidades_correctas ... entidades).
# The edit removes the 'num_credito' row entirely (shifting everything up by
# one row) and replaces the final 'entidades' row's JSON-ish dump with a
# version whose nested dictionaries have their keys reordered.

$ws.Rows.Item(2).Delete()

$newEntidades = '{''ner-general'': {''entidades_correctas'': 0, ''entidades_incorrectas'': 0, ''porcentaje_entidades_correctas'': 1, ''paginas_correctas'': 0, ''paginas_incorrectas'': 0, ''porcentaje_paginas_correctas'': 1}, ''ner-escritura'': {''entidades_correctas'': 34, ''entidades_incorrectas'': 3, ''porcentaje_entidades_correctas'': 0.918918918918919, ''paginas_correctas'': 34, ''paginas_incorrectas'': 3, ''porcentaje_paginas_correctas'': 0.918918918918919, ''DISTRITO_RPP'': {''pagina_encontrado'': 5, ''pagina_es_correcta'': False, ''pagina_real'': 62, ''texto_es_correcto'': False, ''valor_encontrado'': ''CHIHUAHUA,'', ''valor_real'': ''CHIHUAHUA, CHIH''}, ''SUPERFICIE_CONSTRUCCION'': {''pagina_encontrado'': 5, ''pagina_es_correcta'': True, ''pagina_real'': 5, ''texto_es_correcto'': False, ''valor_encontrado'': ''68,711.89 SESENTA Y OCHO MIL SETECIENTOS ONCE METROS OCHENTA Y NUEVE DECIMETROS CUADRADOS'', ''valor_real'': ''120.05 M2. CIENTO VEINTE METROS CINCO DECIMETROS CUADRADOS''}, ''TASA_INTERES_MORATORIO'': {''pagina_encontrado'': 49, ''pagina_es_correcta'': False, ''pagina_real'': 14, ''texto_es_correcto'': False, ''valor_encontrado'': ''9.00% ANUAL FIJO'', ''valor_real'': ''9.00% NUEVE POR CIENTO''}}, ''ner-certificado'': {''entidades_correctas'': 7, ''entidades_incorrectas'': 4, ''porcentaje_entidades_correctas'': 0.6363636363636364, ''paginas_correctas'': 7, ''paginas_incorrectas'': 4, ''porcentaje_paginas_correctas'': 1.0, ''CURP_ACREDITADO'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': False, ''valor_encontrado'': '' '', ''valor_real'': ''OIRF790826HCHLML''}, ''FECHA_EMISION'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': False, ''valor_encontrado'': '' '', ''valor_real'': ''25 DE AGOSTO DE 2017 ''}, ''NUMERO_CREDITO'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': False, ''valor_encontrado'': '' '', ''valor_real'': ''805117974''}, ''RFC_ACREDITADO'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': False, ''valor_encontrado'': ''FACULTADES'', ''valor_real'': ''OIRF790828152''}}, ''avaluo_adjudicacion'': {''entidades_correctas'': 0, ''entidades_incorrectas'': 0, ''porcentaje_entidades_correctas'': 1, ''paginas_correctas'': 0, ''paginas_incorrectas'': 0, ''porcentaje_paginas_correctas'': 1}}'

$ws.Range("B8").Value = $newEntidades
